$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Rename header labels for the squared columns
$ws.Range("D1").Value = "xcua"
$ws.Range("E1").Value = "ycua"
$ws.Range("F1").Value = "xycua"

# Move the active selection to I9, matching the saved cursor position
$ws.Range("I9").Select()
